# Punctul 9 din proiect.docx: "0%" -> "100% - @see Dezvoltare.docx"
# (scris ca trei run-uri separate, fiecare cu w:lang="ro-RO", la fel ca
# structura deja folosita la punctul 10 din document).

$d = $word.ActiveDocument

# Gaseste paragraful punctului 9 printr-un text distinctiv/unic din el.
$anchor = $d.Content
$ok = $anchor.Find.Execute("metodologiei de dezvoltare folosite")
if (-not $ok) {
    throw "Nu s-a gasit textul punctului 9."
}

$anchorPara = $anchor.Paragraphs(1).Range

# In interiorul acelui paragraf, localizeaza exact run-ul "0%".
$target = $d.Range($anchorPara.Start, $anchorPara.End)
$found = $target.Find.Execute("0%")
if (-not $found) {
    throw "Nu s-a gasit '0%' in paragraful punctului 9."
}

# Sterge "0%" si insereaza in locul lui trei run-uri separate:
# "100", "%" si " - @see Dezvoltare.docx", fiecare cu limba ro-RO.
$target.Text = ""
$pos = $target.Start

$run1 = $d.Range($pos, $pos)
$run1.InsertAfter("100")
$run1.LanguageID = "ro-RO"

$pos = $run1.End
$run2 = $d.Range($pos, $pos)
$run2.InsertAfter("%")
$run2.LanguageID = "ro-RO"

$pos = $run2.End
$run3 = $d.Range($pos, $pos)
$run3.InsertAfter(" - @see Dezvoltare.docx")
$run3.LanguageID = "ro-RO"
